# "Finished 60 questions in Excel file for evaluation"
# Fill in the remaining question/answer rows (35-60) for the last six
# research papers, and correct the Question Type on row 33 (index 32).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 (Index 32): Question Type was mislabeled, fix it
$ws.Range("B33").Value = 'Open Ended Question'

# Row 36 (Index 35)
$ws.Range("B36").Value = 'What'
$ws.Range("C36").Value = 'What specific technologies does the smart home system utilize?'
$ws.Range("D36").Value = 'The continuous development of the society has led to the improvement of people''s quality of life and consumption level. At the same time, peoples demand for all aspects of production and life is also increasing, thus promoting the emergence and innovation of intelligent household appliances. To manage these devices conveniently and quickly and enrich family life, "smart home" bureau plays a very important role. Smart home, which enters people''s family life, uses communication technology, Internet connection technology, automatic fire control technology, network wiring technology, and visual and audio transmission technology to communicate with home. Mobile terminals have been developed, and more and more PC functions have been realized. Based on the hardware platform of the smart home management system, two solutions were put forward. The first solution is combined with the current 5G network, and through it, the user can control the mobile phone and other mobile terminals on the corresponding application operation instructions to create. The second solution is the design of the web server intelligent management system, for relevant information. It is collected into the database of the server, allows remote access to the node and subscriber information related to storage in the database through the Internet, and searches the information to control the home lighting and temperature. This system is designed to imitate the modular scheme, which includes the central control module, sensor data acquisition module, and software module. Finally, on the browser side and the electronic devices of Android operating system, it realizes the wireless control of lighting, air conditioning, washing machine, and other devices, as well as the detection of the home environment.'

# Row 37 (Index 36)
$ws.Range("B37").Value = 'Which'
$ws.Range("C37").Value = 'Which solution (5G network or web server) offers remote access for controlling home devices?'
$ws.Range("D37").Value = 'The continuous development of the society has led to the improvement of people''s quality of life and consumption level. At the same time, peoples demand for all aspects of production and life is also increasing, thus promoting the emergence and innovation of intelligent household appliances. To manage these devices conveniently and quickly and enrich family life, "smart home" bureau plays a very important role. Smart home, which enters people''s family life, uses communication technology, Internet connection technology, automatic fire control technology, network wiring technology, and visual and audio transmission technology to communicate with home. Mobile terminals have been developed, and more and more PC functions have been realized. Based on the hardware platform of the smart home management system, two solutions were put forward. The first solution is combined with the current 5G network, and through it, the user can control the mobile phone and other mobile terminals on the corresponding application operation instructions to create. The second solution is the design of the web server intelligent management system, for relevant information. It is collected into the database of the server, allows remote access to the node and subscriber information related to storage in the database through the Internet, and searches the information to control the home lighting and temperature. This system is designed to imitate the modular scheme, which includes the central control module, sensor data acquisition module, and software module. Finally, on the browser side and the electronic devices of Android operating system, it realizes the wireless control of lighting, air conditioning, washing machine, and other devices, as well as the detection of the home environment.'

# Row 38 (Index 37)
$ws.Range("B38").Value = 'When'
$ws.Range("C38").Value = 'When might a user choose to control their smart home devices from a mobile terminal versus a web interface?'
$ws.Range("D38").Value = 'The continuous development of the society has led to the improvement of people''s quality of life and consumption level. At the same time, peoples demand for all aspects of production and life is also increasing, thus promoting the emergence and innovation of intelligent household appliances. To manage these devices conveniently and quickly and enrich family life, "smart home" bureau plays a very important role. Smart home, which enters people''s family life, uses communication technology, Internet connection technology, automatic fire control technology, network wiring technology, and visual and audio transmission technology to communicate with home. Mobile terminals have been developed, and more and more PC functions have been realized. Based on the hardware platform of the smart home management system, two solutions were put forward. The first solution is combined with the current 5G network, and through it, the user can control the mobile phone and other mobile terminals on the corresponding application operation instructions to create. The second solution is the design of the web server intelligent management system, for relevant information. It is collected into the database of the server, allows remote access to the node and subscriber information related to storage in the database through the Internet, and searches the information to control the home lighting and temperature. This system is designed to imitate the modular scheme, which includes the central control module, sensor data acquisition module, and software module. Finally, on the browser side and the electronic devices of Android operating system, it realizes the wireless control of lighting, air conditioning, washing machine, and other devices, as well as the detection of the home environment.'

# Row 39 (Index 38)
$ws.Range("B39").Value = 'Who'
$ws.Range("C39").Value = 'Who are the intended users of this smart home system?'
$ws.Range("D39").Value = 'The continuous development of the society has led to the improvement of people''s quality of life and consumption level. At the same time, peoples demand for all aspects of production and life is also increasing, thus promoting the emergence and innovation of intelligent household appliances. To manage these devices conveniently and quickly and enrich family life, "smart home" bureau plays a very important role. Smart home, which enters people''s family life, uses communication technology, Internet connection technology, automatic fire control technology, network wiring technology, and visual and audio transmission technology to communicate with home. Mobile terminals have been developed, and more and more PC functions have been realized. Based on the hardware platform of the smart home management system, two solutions were put forward. The first solution is combined with the current 5G network, and through it, the user can control the mobile phone and other mobile terminals on the corresponding application operation instructions to create. The second solution is the design of the web server intelligent management system, for relevant information. It is collected into the database of the server, allows remote access to the node and subscriber information related to storage in the database through the Internet, and searches the information to control the home lighting and temperature. This system is designed to imitate the modular scheme, which includes the central control module, sensor data acquisition module, and software module. Finally, on the browser side and the electronic devices of Android operating system, it realizes the wireless control of lighting, air conditioning, washing machine, and other devices, as well as the detection of the home environment.'

# Row 40 (Index 39)
$ws.Range("B40").Value = 'How'
$ws.Range("C40").Value = 'How does the system enable control of lighting, air conditioning, and other appliances?'
$ws.Range("D40").Value = 'The continuous development of the society has led to the improvement of people''s quality of life and consumption level. At the same time, peoples demand for all aspects of production and life is also increasing, thus promoting the emergence and innovation of intelligent household appliances. To manage these devices conveniently and quickly and enrich family life, "smart home" bureau plays a very important role. Smart home, which enters people''s family life, uses communication technology, Internet connection technology, automatic fire control technology, network wiring technology, and visual and audio transmission technology to communicate with home. Mobile terminals have been developed, and more and more PC functions have been realized. Based on the hardware platform of the smart home management system, two solutions were put forward. The first solution is combined with the current 5G network, and through it, the user can control the mobile phone and other mobile terminals on the corresponding application operation instructions to create. The second solution is the design of the web server intelligent management system, for relevant information. It is collected into the database of the server, allows remote access to the node and subscriber information related to storage in the database through the Internet, and searches the information to control the home lighting and temperature. This system is designed to imitate the modular scheme, which includes the central control module, sensor data acquisition module, and software module. Finally, on the browser side and the electronic devices of Android operating system, it realizes the wireless control of lighting, air conditioning, washing machine, and other devices, as well as the detection of the home environment.'

# Row 41 (Index 40)
$ws.Range("B41").Value = 'What'
$ws.Range("C41").Value = 'What is the purpose of integrating space topological cognition into the model?'
$ws.Range("D41").Value = 'Target-driven visual navigation is essential for many applications in robotics, and it has gained increasing interest in recent years. In this work, inspired by animal cognitive mechanisms, we propose a novel navigation architecture that simultaneously learns exploration policy and encodes environmental structure. First, to learn exploration policy directly from raw visual input, we use deep reinforcement learning as the basic framework and allow agents to create rewards for themselves as learning signals. In our approach, the reward for the current observation is driven by curiosity and calculated by a count-based approach and temporal distance. While agents learn exploration policy, we use temporal distance to find waypoints in observation sequences and incrementally describe the structure of the environment in a way that integrates episodic memory. Finally, space topological cognition is integrated into the model as a path planning module and combined with a locomotion network to obtain a more generalized approach to navigation. We test our approach in the DMlab, a visually rich 3D environment, and validate its exploration efficiency and navigation performance through extensive experiments. The experimental results show that our approach can explore and encode the environment more efficiently and has better capability in dealing with stochastic objects. In navigation tasks, agents can use space topological cognition to effectively reach the target and guide detour behaviour when a path is unavailable, exhibiting good environmental adaptability.'

# Row 42 (Index 41)
$ws.Range("B42").Value = 'Which'
$ws.Range("C42").Value = 'Which component of the architecture is responsible for guiding detour behavior?'
$ws.Range("D42").Value = 'Target-driven visual navigation is essential for many applications in robotics, and it has gained increasing interest in recent years. In this work, inspired by animal cognitive mechanisms, we propose a novel navigation architecture that simultaneously learns exploration policy and encodes environmental structure. First, to learn exploration policy directly from raw visual input, we use deep reinforcement learning as the basic framework and allow agents to create rewards for themselves as learning signals. In our approach, the reward for the current observation is driven by curiosity and calculated by a count-based approach and temporal distance. While agents learn exploration policy, we use temporal distance to find waypoints in observation sequences and incrementally describe the structure of the environment in a way that integrates episodic memory. Finally, space topological cognition is integrated into the model as a path planning module and combined with a locomotion network to obtain a more generalized approach to navigation. We test our approach in the DMlab, a visually rich 3D environment, and validate its exploration efficiency and navigation performance through extensive experiments. The experimental results show that our approach can explore and encode the environment more efficiently and has better capability in dealing with stochastic objects. In navigation tasks, agents can use space topological cognition to effectively reach the target and guide detour behaviour when a path is unavailable, exhibiting good environmental adaptability.'

# Row 43 (Index 42)
$ws.Range("B43").Value = 'When'
$ws.Range("C43").Value = 'When would an agent rely on space topological cognition instead of direct locomotion for navigation?'
$ws.Range("D43").Value = 'Target-driven visual navigation is essential for many applications in robotics, and it has gained increasing interest in recent years. In this work, inspired by animal cognitive mechanisms, we propose a novel navigation architecture that simultaneously learns exploration policy and encodes environmental structure. First, to learn exploration policy directly from raw visual input, we use deep reinforcement learning as the basic framework and allow agents to create rewards for themselves as learning signals. In our approach, the reward for the current observation is driven by curiosity and calculated by a count-based approach and temporal distance. While agents learn exploration policy, we use temporal distance to find waypoints in observation sequences and incrementally describe the structure of the environment in a way that integrates episodic memory. Finally, space topological cognition is integrated into the model as a path planning module and combined with a locomotion network to obtain a more generalized approach to navigation. We test our approach in the DMlab, a visually rich 3D environment, and validate its exploration efficiency and navigation performance through extensive experiments. The experimental results show that our approach can explore and encode the environment more efficiently and has better capability in dealing with stochastic objects. In navigation tasks, agents can use space topological cognition to effectively reach the target and guide detour behaviour when a path is unavailable, exhibiting good environmental adaptability.'

# Row 44 (Index 43)
$ws.Range("B44").Value = 'How'
$ws.Range("C44").Value = 'How do the curiosity-driven rewards influence the agent''s exploration behavior?'
$ws.Range("D44").Value = 'Target-driven visual navigation is essential for many applications in robotics, and it has gained increasing interest in recent years. In this work, inspired by animal cognitive mechanisms, we propose a novel navigation architecture that simultaneously learns exploration policy and encodes environmental structure. First, to learn exploration policy directly from raw visual input, we use deep reinforcement learning as the basic framework and allow agents to create rewards for themselves as learning signals. In our approach, the reward for the current observation is driven by curiosity and calculated by a count-based approach and temporal distance. While agents learn exploration policy, we use temporal distance to find waypoints in observation sequences and incrementally describe the structure of the environment in a way that integrates episodic memory. Finally, space topological cognition is integrated into the model as a path planning module and combined with a locomotion network to obtain a more generalized approach to navigation. We test our approach in the DMlab, a visually rich 3D environment, and validate its exploration efficiency and navigation performance through extensive experiments. The experimental results show that our approach can explore and encode the environment more efficiently and has better capability in dealing with stochastic objects. In navigation tasks, agents can use space topological cognition to effectively reach the target and guide detour behaviour when a path is unavailable, exhibiting good environmental adaptability.'

# Row 45 (Index 44)
$ws.Range("B45").Value = 'Hypothetical'
$ws.Range("C45").Value = 'What emerging field of study addresses the limitations of traditional safety ergonomics in the era of big data, aiming to enhance both safety and cleaner production practices?'
$ws.Range("D45").Value = 'Safety ergonomics is an important branch of safety science and environmental engineering. As humans enter the era of big data, the development of information technology has brought new opportunities and challenges to the innovation, transformation, and upgrading of safety ergonomics, as the traditional safety ergonomics theory has gradually failed to adapt to the need for safe and clean production. Intelligent safety ergonomics (ISE) is regarded as a new direction for the development of safety ergonomics in the era of big data. Unfortunately, since ISE is an emerging concept, there is no research to clarify its basic problems, which leads to a lack of theoretical guidance for the research and practice of ISE. In order to solve the shortcomings of traditional safety ergonomics theories and methods, first of all, this paper answers the basic questions of ISE, including the basic concepts, characteristics, attributes, contents, and research objects. Then, practical application functions of ISE are systematically clarified. Finally, following the life cycle of the design, implementation, operation, and maintenance of the system, it ends with a discussion of the challenges and application prospects of ISE. The conclusion shows that ISE is a cleaner research direction for ergonomics in the era of big data, that it can deepen the understanding of humans, machines, and environment systems, and it can provide a new method for further research on safety and cleaner production. Overall, this paper not only helps safety researchers and practitioners to correctly understand the concept of intelligent safety ergonomics, but it will certainly inject energy and vitality into the development of safety ergonomics and cleaner production.'

# Row 46 (Index 45)
$ws.Range("B46").Value = 'Hypothetical'
$ws.Range("C46").Value = 'How can the integration of information technology and safety science lead to a more holistic understanding of human-machine-environment systems, ultimately improving workplace safety and sustainability?'
$ws.Range("D46").Value = 'Safety ergonomics is an important branch of safety science and environmental engineering. As humans enter the era of big data, the development of information technology has brought new opportunities and challenges to the innovation, transformation, and upgrading of safety ergonomics, as the traditional safety ergonomics theory has gradually failed to adapt to the need for safe and clean production. Intelligent safety ergonomics (ISE) is regarded as a new direction for the development of safety ergonomics in the era of big data. Unfortunately, since ISE is an emerging concept, there is no research to clarify its basic problems, which leads to a lack of theoretical guidance for the research and practice of ISE. In order to solve the shortcomings of traditional safety ergonomics theories and methods, first of all, this paper answers the basic questions of ISE, including the basic concepts, characteristics, attributes, contents, and research objects. Then, practical application functions of ISE are systematically clarified. Finally, following the life cycle of the design, implementation, operation, and maintenance of the system, it ends with a discussion of the challenges and application prospects of ISE. The conclusion shows that ISE is a cleaner research direction for ergonomics in the era of big data, that it can deepen the understanding of humans, machines, and environment systems, and it can provide a new method for further research on safety and cleaner production. Overall, this paper not only helps safety researchers and practitioners to correctly understand the concept of intelligent safety ergonomics, but it will certainly inject energy and vitality into the development of safety ergonomics and cleaner production.'

# Row 47 (Index 46)
$ws.Range("B47").Value = 'Open Ended Question'
$ws.Range("C47").Value = 'How does the study challenge previous findings about the role of goal-driven manipulation in location-reward association?'
$ws.Range("D47").Value = 'Many studies have reported attentional biases based on feature-reward associations. However, the effects of location-reward associations on attentional selection remain less well-understood. Unlike feature cases, a previous study that induced participants'' awareness of the location-reward association by instructing them to look for a high-reward location has suggested the critical role of goal-driven manipulations in such associations. In this study, we investigated whether the reward effect occurred without goal-driven manipulations if participants were spontaneously aware of the location-reward association. We conducted three experiments using a visual search task that included four circles where participants received rewards; one possible target location was associated with a high reward, and another with a low reward. In Experiment 1, the target was presented among distractors, and participants had to search for the target. The results showed a faster reaction time in the high-reward rather than the low-reward locations only in participants aware of the location-reward association, even if they were not required to look for the association. Moreover, in Experiment 2, we replicated the main findings of Experiment 1, even when the target had an abrupt visual onset to restrict goal-driven manipulations. Furthermore, Experiment 3 confirmed that the effect observed in Experiment 2 could not be attributed to the initial eye position. These findings suggest that goal-driven manipulations are unnecessary for inducing reward biases to high-reward locations. We concluded that awareness of the association rather than goal-driven manipulations is crucial for the location-reward effect.'

# Row 48 (Index 47)
$ws.Range("B48").Value = 'Yes / No question'
$ws.Range("C48").Value = 'Did the researchers find a significant location-reward effect even when participants were not actively searching for high-reward locations?'
$ws.Range("D48").Value = 'Many studies have reported attentional biases based on feature-reward associations. However, the effects of location-reward associations on attentional selection remain less well-understood. Unlike feature cases, a previous study that induced participants'' awareness of the location-reward association by instructing them to look for a high-reward location has suggested the critical role of goal-driven manipulations in such associations. In this study, we investigated whether the reward effect occurred without goal-driven manipulations if participants were spontaneously aware of the location-reward association. We conducted three experiments using a visual search task that included four circles where participants received rewards; one possible target location was associated with a high reward, and another with a low reward. In Experiment 1, the target was presented among distractors, and participants had to search for the target. The results showed a faster reaction time in the high-reward rather than the low-reward locations only in participants aware of the location-reward association, even if they were not required to look for the association. Moreover, in Experiment 2, we replicated the main findings of Experiment 1, even when the target had an abrupt visual onset to restrict goal-driven manipulations. Furthermore, Experiment 3 confirmed that the effect observed in Experiment 2 could not be attributed to the initial eye position. These findings suggest that goal-driven manipulations are unnecessary for inducing reward biases to high-reward locations. We concluded that awareness of the association rather than goal-driven manipulations is crucial for the location-reward effect.'

# Row 49 (Index 48)
$ws.Range("B49").Value = 'Multiple choice question'
$ws.Range("C49").Value = 'Which factor appears to be the most important for inducing a reward bias in the study?
a) Presence of a reward
b) Awareness of location-reward association
c) Goal-driven search for high-reward locations
d) Abrupt onset of the target'
$ws.Range("D49").Value = 'Many studies have reported attentional biases based on feature-reward associations. However, the effects of location-reward associations on attentional selection remain less well-understood. Unlike feature cases, a previous study that induced participants'' awareness of the location-reward association by instructing them to look for a high-reward location has suggested the critical role of goal-driven manipulations in such associations. In this study, we investigated whether the reward effect occurred without goal-driven manipulations if participants were spontaneously aware of the location-reward association. We conducted three experiments using a visual search task that included four circles where participants received rewards; one possible target location was associated with a high reward, and another with a low reward. In Experiment 1, the target was presented among distractors, and participants had to search for the target. The results showed a faster reaction time in the high-reward rather than the low-reward locations only in participants aware of the location-reward association, even if they were not required to look for the association. Moreover, in Experiment 2, we replicated the main findings of Experiment 1, even when the target had an abrupt visual onset to restrict goal-driven manipulations. Furthermore, Experiment 3 confirmed that the effect observed in Experiment 2 could not be attributed to the initial eye position. These findings suggest that goal-driven manipulations are unnecessary for inducing reward biases to high-reward locations. We concluded that awareness of the association rather than goal-driven manipulations is crucial for the location-reward effect.'

# Row 50 (Index 49)
$ws.Range("B50").Value = 'What'
$ws.Range("C50").Value = 'What were the primary goals of the Hungarian Primary Care Model Programme?'
$ws.Range("D50").Value = 'Background: A Primary Care Model Programme was implemented in Hungary between 2013 and 2017 in order to increase access of disadvantaged population groups to primary care and to offer new preventive services for all clients. In a country with single-handed practices, four group practices or GP clusters were created in the Programme. Six GPs comprised one cluster who together employed nonmedical health professionals and nonprofessional health mediators, the latter recruited from the serviced communities, many of them of Roma ethnicity. Health mediators were tasked by improving access of the local communities - including its vulnerable Roma members - to existing and new services. Health mediators were interviewed about their work experiences, motivation, and overall opinion as members of the clusters as part of the Programme evaluation.
Methods: As part of the Programme evaluation, structured interviews were conducted with all 40 health mediators employed at the time in the Programme. Interviews were transcribed and content analysis was carried out.
Results: Three themes emerged from the transcripts. The first focused on the health mediators'' personal characteristics such as motivation to join the Programme, the way their job increased their self-esteem, social status and health consciousness. Domains of the second theme of their work included importance of on-the-job training and of their insider knowledge of local communities, as well as their pride to have become members of the primary care team. The third theme covered overall functioning of the Programme of which they had mostly positive opinions, notwithstanding some criticism regarding procurement.
Conclusions: Health mediators had earlier worked in various European countries specifically to improve access of Roma ethnic groups to health services but the Hungarian Model Programme was globally the first in which health mediators as non-professional workers became equal members of the primary care team as employees. Their contribution and overwhelmingly positive experiences, along with their useful insights for improvement call for the establishment and funding of health mediator positions in primary care especially in areas with large numbers of disadvantaged Roma populations.'

# Row 51 (Index 50)
$ws.Range("B51").Value = 'Who'
$ws.Range("C51").Value = 'Who were the primary beneficiaries of the health mediator''s work?'
$ws.Range("D51").Value = 'Background: A Primary Care Model Programme was implemented in Hungary between 2013 and 2017 in order to increase access of disadvantaged population groups to primary care and to offer new preventive services for all clients. In a country with single-handed practices, four group practices or GP clusters were created in the Programme. Six GPs comprised one cluster who together employed nonmedical health professionals and nonprofessional health mediators, the latter recruited from the serviced communities, many of them of Roma ethnicity. Health mediators were tasked by improving access of the local communities - including its vulnerable Roma members - to existing and new services. Health mediators were interviewed about their work experiences, motivation, and overall opinion as members of the clusters as part of the Programme evaluation.
Methods: As part of the Programme evaluation, structured interviews were conducted with all 40 health mediators employed at the time in the Programme. Interviews were transcribed and content analysis was carried out.
Results: Three themes emerged from the transcripts. The first focused on the health mediators'' personal characteristics such as motivation to join the Programme, the way their job increased their self-esteem, social status and health consciousness. Domains of the second theme of their work included importance of on-the-job training and of their insider knowledge of local communities, as well as their pride to have become members of the primary care team. The third theme covered overall functioning of the Programme of which they had mostly positive opinions, notwithstanding some criticism regarding procurement.
Conclusions: Health mediators had earlier worked in various European countries specifically to improve access of Roma ethnic groups to health services but the Hungarian Model Programme was globally the first in which health mediators as non-professional workers became equal members of the primary care team as employees. Their contribution and overwhelmingly positive experiences, along with their useful insights for improvement call for the establishment and funding of health mediator positions in primary care especially in areas with large numbers of disadvantaged Roma populations.'

# Row 52 (Index 51)
$ws.Range("B52").Value = 'How'
$ws.Range("C52").Value = 'How did the health mediators'' work increase their self-esteem and social status?'
$ws.Range("D52").Value = 'Background: A Primary Care Model Programme was implemented in Hungary between 2013 and 2017 in order to increase access of disadvantaged population groups to primary care and to offer new preventive services for all clients. In a country with single-handed practices, four group practices or GP clusters were created in the Programme. Six GPs comprised one cluster who together employed nonmedical health professionals and nonprofessional health mediators, the latter recruited from the serviced communities, many of them of Roma ethnicity. Health mediators were tasked by improving access of the local communities - including its vulnerable Roma members - to existing and new services. Health mediators were interviewed about their work experiences, motivation, and overall opinion as members of the clusters as part of the Programme evaluation.
Methods: As part of the Programme evaluation, structured interviews were conducted with all 40 health mediators employed at the time in the Programme. Interviews were transcribed and content analysis was carried out.
Results: Three themes emerged from the transcripts. The first focused on the health mediators'' personal characteristics such as motivation to join the Programme, the way their job increased their self-esteem, social status and health consciousness. Domains of the second theme of their work included importance of on-the-job training and of their insider knowledge of local communities, as well as their pride to have become members of the primary care team. The third theme covered overall functioning of the Programme of which they had mostly positive opinions, notwithstanding some criticism regarding procurement.
Conclusions: Health mediators had earlier worked in various European countries specifically to improve access of Roma ethnic groups to health services but the Hungarian Model Programme was globally the first in which health mediators as non-professional workers became equal members of the primary care team as employees. Their contribution and overwhelmingly positive experiences, along with their useful insights for improvement call for the establishment and funding of health mediator positions in primary care especially in areas with large numbers of disadvantaged Roma populations.'

# Row 53 (Index 52)
$ws.Range("B53").Value = 'Why'
$ws.Range("C53").Value = 'Why is the use of health mediators particularly important in areas with disadvantaged Roma populations?'
$ws.Range("D53").Value = 'Background: A Primary Care Model Programme was implemented in Hungary between 2013 and 2017 in order to increase access of disadvantaged population groups to primary care and to offer new preventive services for all clients. In a country with single-handed practices, four group practices or GP clusters were created in the Programme. Six GPs comprised one cluster who together employed nonmedical health professionals and nonprofessional health mediators, the latter recruited from the serviced communities, many of them of Roma ethnicity. Health mediators were tasked by improving access of the local communities - including its vulnerable Roma members - to existing and new services. Health mediators were interviewed about their work experiences, motivation, and overall opinion as members of the clusters as part of the Programme evaluation.
Methods: As part of the Programme evaluation, structured interviews were conducted with all 40 health mediators employed at the time in the Programme. Interviews were transcribed and content analysis was carried out.
Results: Three themes emerged from the transcripts. The first focused on the health mediators'' personal characteristics such as motivation to join the Programme, the way their job increased their self-esteem, social status and health consciousness. Domains of the second theme of their work included importance of on-the-job training and of their insider knowledge of local communities, as well as their pride to have become members of the primary care team. The third theme covered overall functioning of the Programme of which they had mostly positive opinions, notwithstanding some criticism regarding procurement.
Conclusions: Health mediators had earlier worked in various European countries specifically to improve access of Roma ethnic groups to health services but the Hungarian Model Programme was globally the first in which health mediators as non-professional workers became equal members of the primary care team as employees. Their contribution and overwhelmingly positive experiences, along with their useful insights for improvement call for the establishment and funding of health mediator positions in primary care especially in areas with large numbers of disadvantaged Roma populations.'

# Row 54 (Index 53)
$ws.Range("B54").Value = 'Where'
$ws.Range("C54").Value = 'Where were the GP clusters implementing the Primary Care Model Programme located?'
$ws.Range("D54").Value = 'Background: A Primary Care Model Programme was implemented in Hungary between 2013 and 2017 in order to increase access of disadvantaged population groups to primary care and to offer new preventive services for all clients. In a country with single-handed practices, four group practices or GP clusters were created in the Programme. Six GPs comprised one cluster who together employed nonmedical health professionals and nonprofessional health mediators, the latter recruited from the serviced communities, many of them of Roma ethnicity. Health mediators were tasked by improving access of the local communities - including its vulnerable Roma members - to existing and new services. Health mediators were interviewed about their work experiences, motivation, and overall opinion as members of the clusters as part of the Programme evaluation.
Methods: As part of the Programme evaluation, structured interviews were conducted with all 40 health mediators employed at the time in the Programme. Interviews were transcribed and content analysis was carried out.
Results: Three themes emerged from the transcripts. The first focused on the health mediators'' personal characteristics such as motivation to join the Programme, the way their job increased their self-esteem, social status and health consciousness. Domains of the second theme of their work included importance of on-the-job training and of their insider knowledge of local communities, as well as their pride to have become members of the primary care team. The third theme covered overall functioning of the Programme of which they had mostly positive opinions, notwithstanding some criticism regarding procurement.
Conclusions: Health mediators had earlier worked in various European countries specifically to improve access of Roma ethnic groups to health services but the Hungarian Model Programme was globally the first in which health mediators as non-professional workers became equal members of the primary care team as employees. Their contribution and overwhelmingly positive experiences, along with their useful insights for improvement call for the establishment and funding of health mediator positions in primary care especially in areas with large numbers of disadvantaged Roma populations.'

# Row 55 (Index 54)
$ws.Range("B55").Value = 'When'
$ws.Range("C55").Value = 'When did the evaluation of the Primary Care Model Programme take place?'
$ws.Range("D55").Value = 'Background: A Primary Care Model Programme was implemented in Hungary between 2013 and 2017 in order to increase access of disadvantaged population groups to primary care and to offer new preventive services for all clients. In a country with single-handed practices, four group practices or GP clusters were created in the Programme. Six GPs comprised one cluster who together employed nonmedical health professionals and nonprofessional health mediators, the latter recruited from the serviced communities, many of them of Roma ethnicity. Health mediators were tasked by improving access of the local communities - including its vulnerable Roma members - to existing and new services. Health mediators were interviewed about their work experiences, motivation, and overall opinion as members of the clusters as part of the Programme evaluation.
Methods: As part of the Programme evaluation, structured interviews were conducted with all 40 health mediators employed at the time in the Programme. Interviews were transcribed and content analysis was carried out.
Results: Three themes emerged from the transcripts. The first focused on the health mediators'' personal characteristics such as motivation to join the Programme, the way their job increased their self-esteem, social status and health consciousness. Domains of the second theme of their work included importance of on-the-job training and of their insider knowledge of local communities, as well as their pride to have become members of the primary care team. The third theme covered overall functioning of the Programme of which they had mostly positive opinions, notwithstanding some criticism regarding procurement.
Conclusions: Health mediators had earlier worked in various European countries specifically to improve access of Roma ethnic groups to health services but the Hungarian Model Programme was globally the first in which health mediators as non-professional workers became equal members of the primary care team as employees. Their contribution and overwhelmingly positive experiences, along with their useful insights for improvement call for the establishment and funding of health mediator positions in primary care especially in areas with large numbers of disadvantaged Roma populations.'

# Row 56 (Index 55)
$ws.Range("B56").Value = 'Which'
$ws.Range("C56").Value = 'Which of the health mediators'' job aspects received the most positive feedback?'
$ws.Range("D56").Value = 'Background: A Primary Care Model Programme was implemented in Hungary between 2013 and 2017 in order to increase access of disadvantaged population groups to primary care and to offer new preventive services for all clients. In a country with single-handed practices, four group practices or GP clusters were created in the Programme. Six GPs comprised one cluster who together employed nonmedical health professionals and nonprofessional health mediators, the latter recruited from the serviced communities, many of them of Roma ethnicity. Health mediators were tasked by improving access of the local communities - including its vulnerable Roma members - to existing and new services. Health mediators were interviewed about their work experiences, motivation, and overall opinion as members of the clusters as part of the Programme evaluation.
Methods: As part of the Programme evaluation, structured interviews were conducted with all 40 health mediators employed at the time in the Programme. Interviews were transcribed and content analysis was carried out.
Results: Three themes emerged from the transcripts. The first focused on the health mediators'' personal characteristics such as motivation to join the Programme, the way their job increased their self-esteem, social status and health consciousness. Domains of the second theme of their work included importance of on-the-job training and of their insider knowledge of local communities, as well as their pride to have become members of the primary care team. The third theme covered overall functioning of the Programme of which they had mostly positive opinions, notwithstanding some criticism regarding procurement.
Conclusions: Health mediators had earlier worked in various European countries specifically to improve access of Roma ethnic groups to health services but the Hungarian Model Programme was globally the first in which health mediators as non-professional workers became equal members of the primary care team as employees. Their contribution and overwhelmingly positive experiences, along with their useful insights for improvement call for the establishment and funding of health mediator positions in primary care especially in areas with large numbers of disadvantaged Roma populations.'

# Row 57 (Index 56)
$ws.Range("B57").Value = 'Open Ended Question'
$ws.Range("C57").Value = 'How can clinicians improve their communication with patients to cultivate stronger relationships and promote better adherence to dietary recommendations? '
$ws.Range("D57").Value = 'For clinicians, effective communication goes beyond the delivery of scientific information to include an understanding of who the patient is and what they value; treating them with respect; and acknowledging their emotional and social realities. Recognizing our own implicit biases and cultivating a more mindful approach to the impact of language, especially around weight and food, and using verbal and nonverbal approaches to convey empathic concern can improve our relationships with our patients. Patients who feel seen and respected are more likely to follow through with recommendations including dietary change, leading to improved health and quality of life.'

# Row 58 (Index 57)
$ws.Range("B58").Value = 'Why'
$ws.Range("C58").Value = 'Why is it crucial for clinicians to recognize their own implicit biases about weight and food?'
$ws.Range("D58").Value = 'For clinicians, effective communication goes beyond the delivery of scientific information to include an understanding of who the patient is and what they value; treating them with respect; and acknowledging their emotional and social realities. Recognizing our own implicit biases and cultivating a more mindful approach to the impact of language, especially around weight and food, and using verbal and nonverbal approaches to convey empathic concern can improve our relationships with our patients. Patients who feel seen and respected are more likely to follow through with recommendations including dietary change, leading to improved health and quality of life.'

# Row 59 (Index 58)
$ws.Range("B59").Value = 'Open Ended Question'
$ws.Range("C59").Value = 'How does a leader''s level of trait empathy influence their emotional reactions and perceived effectiveness after providing negative feedback?'
$ws.Range("D59").Value = 'Although providing negative performance feedback can enhance employee performance, leaders are sometimes reluctant to engage in this activity. Reflecting this, prior research has identified negative feedback provision as an aversive, yet potentially rewarding, managerial activity. However, little is known about how providing negative feedback impacts the effectiveness of leaders who do so. To shed light on this issue, we develop and test a theoretical model that identifies how leaders'' proximal and distal reactions to providing negative feedback are contingent upon their levels of trait empathy. Supporting our theory, results from an experience sampling study indicate that leaders higher in trait empathy report feeling both less attentive and more distressed after providing subordinates with negative feedback, whereas leaders lower in trait empathy report feeling more attentive and less distressed. Attentiveness and distress, in turn, were associated with leaders'' daily perceptions of their effectiveness; distress was also associated with leaders'' daily enactment of transformational leadership behavior. Results of two subsequent studies focused on single episodes of negative feedback provision revealed that trait empathy amplifies the extent to which feedback recipients'' negative emotional reactions impact additional leader effectiveness criteria (e.g., executive functioning and planning/problem-solving), further supporting the need to account for the crucial role of trait empathy in the feedback-provision process. Altogether, our research provides a novel perspective on the feedback-giving process by shifting the focus of theorizing from the recipient to the provider, while challenging current thinking about leader empathy by highlighting its potential downside for leadership. (PsycInfo Database Record (c) 2022 APA, all rights reserved).'

# Row 60 (Index 59)
$ws.Range("B60").Value = 'Yes / No question'
$ws.Range("C60").Value = 'Does the study suggest that there are potential downsides to leader empathy in the context of giving negative feedback?'
$ws.Range("D60").Value = 'Although providing negative performance feedback can enhance employee performance, leaders are sometimes reluctant to engage in this activity. Reflecting this, prior research has identified negative feedback provision as an aversive, yet potentially rewarding, managerial activity. However, little is known about how providing negative feedback impacts the effectiveness of leaders who do so. To shed light on this issue, we develop and test a theoretical model that identifies how leaders'' proximal and distal reactions to providing negative feedback are contingent upon their levels of trait empathy. Supporting our theory, results from an experience sampling study indicate that leaders higher in trait empathy report feeling both less attentive and more distressed after providing subordinates with negative feedback, whereas leaders lower in trait empathy report feeling more attentive and less distressed. Attentiveness and distress, in turn, were associated with leaders'' daily perceptions of their effectiveness; distress was also associated with leaders'' daily enactment of transformational leadership behavior. Results of two subsequent studies focused on single episodes of negative feedback provision revealed that trait empathy amplifies the extent to which feedback recipients'' negative emotional reactions impact additional leader effectiveness criteria (e.g., executive functioning and planning/problem-solving), further supporting the need to account for the crucial role of trait empathy in the feedback-provision process. Altogether, our research provides a novel perspective on the feedback-giving process by shifting the focus of theorizing from the recipient to the provider, while challenging current thinking about leader empathy by highlighting its potential downside for leadership. (PsycInfo Database Record (c) 2022 APA, all rights reserved).'

# Row 61 (Index 60)
$ws.Range("B61").Value = 'Open Ended Question'
$ws.Range("C61").Value = 'How can nurse leaders use emotional intelligence to specifically address the challenges of stress, exhaustion, and the risk of moral injury brought on by the COVID-19 pandemic?'
$ws.Range("D61").Value = 'Emotionally intelligent leaders demonstrate a sensitivity to their own and other people''s psychological health and well-being, directing others towards common goals while developing effective personal relationships with their colleagues and team members. Emotional intelligence is particularly relevant in the context of the coronavirus disease 2019 pandemic, where nurse leaders need to demonstrate this skill when supporting their teams to manage high levels of stress, exhaustion and the risk of moral injury. This article explores emotional intelligence, discusses its importance as a characteristic of effective nurse leaders and managers, and suggests practical activities that leaders can undertake to develop their emotional intelligence skills.'

# Move the cursor/scroll position back to the top of the sheet, since the
# previous in-progress selection (D35 / scrolled to row 23) is now stale.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A1").Select()

Write-Host "Applied 60-question dataset updates."